$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B15").Value = "F1"
$ws.Range("A15").Value = "500mA"
$ws.Range("C15").Value = "Fuse_0805_2012Metric_Pad1.15x1.40mm_HandSolder"
$ws.Range("D15").Value = "C2649565"

$ws.Range("A1:E15").RowHeight = 20.1

$ws.Range("D15").Select()
